# The deck ships two DrawingML theme parts:
#   ppt/theme/theme1.xml -> "Office Theme" / clrScheme "Office"      (currently only wired to the Notes Master)
#   ppt/theme/theme2.xml -> "Integral"     / clrScheme "Red Violet"  (the theme actually driving the slide master / slides)
#
# The authored edit swaps the two themes' content, so the slides end up styled with the
# plain "Office" colour palette instead of the pink/purple "Red Violet" one (and the Notes
# Master ends up with the "Red Violet" palette). The font scheme and format (fill/line/effect)
# scheme are identical between the two theme parts already, so the only real difference is the
# 12-slot colour scheme (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink) - that's what we rewrite
# on the presentation's active theme via the PowerPoint colour-scheme object model.

$p = $ppt.ActivePresentation
$cs = $p.ColorSchemes.Item(1)

# Target palette = the current theme1.xml ("Office Theme" / clrScheme "Office") values,
# expressed as VBA-style RGB() integers (R + G*256 + B*65536).
$cs.Colors(1).RGB  = 0x00 + 0x00*256 + 0x00*65536    # dk1      000000
$cs.Colors(2).RGB  = 0xFF + 0xFF*256 + 0xFF*65536    # lt1      FFFFFF
$cs.Colors(3).RGB  = 0x44 + 0x54*256 + 0x6A*65536    # dk2      44546A
$cs.Colors(4).RGB  = 0xE7 + 0xE6*256 + 0xE6*65536    # lt2      E7E6E6
$cs.Colors(5).RGB  = 0x5B + 0x9B*256 + 0xD5*65536    # accent1  5B9BD5
$cs.Colors(6).RGB  = 0xED + 0x7D*256 + 0x31*65536    # accent2  ED7D31
$cs.Colors(7).RGB  = 0xA5 + 0xA5*256 + 0xA5*65536    # accent3  A5A5A5
$cs.Colors(8).RGB  = 0xFF + 0xC0*256 + 0x00*65536    # accent4  FFC000
$cs.Colors(9).RGB  = 0x44 + 0x72*256 + 0xC4*65536    # accent5  4472C4
$cs.Colors(10).RGB = 0x70 + 0xAD*256 + 0x47*65536    # accent6  70AD47
$cs.Colors(11).RGB = 0x05 + 0x63*256 + 0xC1*65536    # hlink    0563C1
$cs.Colors(12).RGB = 0x95 + 0x4F*256 + 0x72*65536    # folHlink 954F72
